$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 5
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 1

$ws.Range("A3").Select()
